$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.777.36"
$ws.Range("E2").Value = "  -4.04%  "
$ws.Range("D3").Value = "2.331.09"
$ws.Range("E3").Value = "  -4.78%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "84.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.525"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.40%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0814"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.46%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").Value = "2.689.69"
$ws.Range("E13").Value = "  -4.77%  "
$ws.Range("E14").Value = "  -5.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("D16").Value = "2.335.80"
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.752"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "39.758.45"
$ws.Range("E18").Value = "  -3.80%  "
$ws.Range("E19").Value = "  -3.10%  "
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  -6.87%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "151.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.36%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.03%  "
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0713"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.69%  "
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0982"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.84%  "
$ws.Range("E38").Value = "  -7.30%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.79%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.17%  "
$ws.Range("D42").Value = "1.931.45"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.02%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.28%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0260"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.17%  "
$ws.Range("E47").Value = "  -7.90%  "
$ws.Range("D48").Value = "2.548.84"
$ws.Range("E48").Value = "  -4.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "91.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.15%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.96%  "
